$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Progress values (column C) updates ---
$ws.Range("C5").Value = 0.95
$ws.Range("C6").Value = 1
$ws.Range("C8").Value = 0.9
$ws.Range("C9").Value = 0.9
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 0.6
$ws.Range("C14").Value = 0.8
$ws.Range("C15").Value = 0.9
$ws.Range("C16").Value = 1
$ws.Range("C18").Value = 0.9
$ws.Range("C19").Value = 1
$ws.Range("C20").Value = 0.3
$ws.Range("C21").Value = 0.5
$ws.Range("C22").Value = 1
$ws.Range("C24").Value = 0.6
$ws.Range("C25").Value = 0.1
$ws.Range("C26").Value = 0.1

# --- Comments (column D) updates ---
# D6 comment removed (was about labelEncoder, no longer relevant)
$ws.Range("D6").Value = ""

# D14 comment rewritten (pipeline remark)
$ws.Range("D14").Value = "Implémentation dans un pipline et peut-être revoir pour l'améliorer"

# D8 comment rewritten
$ws.Range("D8").Value = "Revoir la méthode peux être améliorer / pas nécessaire sur le random forest"

# D9 new comment
$ws.Range("D9").Value = "Oui mais perfectible je pense"

# D18 new comment
$ws.Range("D18").Value = "Fbeta avec beta = 2"

# D24 new comment
$ws.Range("D24").Value = "max_depth principalement"

# D20 new comment
$ws.Range("D20").Value = "Faire un test avec un CV ?"

# --- Row height for row 8 ---
$ws.Rows(8).RowHeight = 56

# --- Sheet view: scroll/selection update ---
$ws.Range("F20").Select()
